# Add carjacking data for 2022-08-22 (updates the "through August 13" running
# report to "through August 14", incrementing one incident in each affected
# neighborhood/month cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the running-report date in the sheet name and the column header text.
$ws.Name = "Through 2022-08-14"
$ws.Range("B1").Value = "August 2022 (through August 14)"

# Cells whose counts increment by one existing incident.
$increments = @{
    "AH2"  = 2
    "AH4"  = 2
    "B5"   = 5
    "B7"   = 6
    "J7"   = 3
    "AP12" = 3
    "AX29" = 2
    "B39"  = 2
    "R61"  = 2
}
foreach ($addr in $increments.Keys) {
    $ws.Range($addr).Value = $increments[$addr]
}

# Cells that are newly populated with a count of 1.
$newCells = @(
    "AX5",
    "AX8",
    "B9",
    "Z9",
    "J15",
    "R16",
    "Z20",
    "BF33",
    "AP46",
    "AP66",
    "AX75",
    "Z76",
    "B96"
)
foreach ($addr in $newCells) {
    $ws.Range($addr).Value = 1
}
